# 6.7 Add Investigate Panel
# Fix misspelled "Avatar" values in column C (Quan-Regualr / Quna-Regular -> Quan-Regular)
# and move the active selection to C18 to reflect the newly added investigate rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Quan-Regular"
$ws.Range("C8").Value = "Quan-Regular"
$ws.Range("C9").Value = "Quan-Regular"
$ws.Range("C14").Value = "Quan-Regular"
$ws.Range("C16").Value = "Quan-Regular"

$ws.Range("C18").Select()
